$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# cryptos list refresh (GitHub Actions data pull).
# Column D holds prices as plain text (e.g. "1.00", "24.31"), not
# numbers -- mark the touched Price cells as Text first so Excel
# does not silently reinterpret numeric-looking strings as numbers.

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '60.164.89'
$ws.Range("E2").Value = '  -6.20%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.276.60'
$ws.Range("E3").Value = '  -5.78%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.18%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '557.92'
$ws.Range("E5").Value = '  -4.60%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '128.21'
$ws.Range("E6").Value = '  -2.51%  '

# Row 7
$ws.Range("E7").Value = '  +0.12%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.268.10'
$ws.Range("E8").Value = '  -6.00%  '

# Row 9
$ws.Range("E9").Value = '  -1.80%  '

# Row 10
$ws.Range("E10").Value = '  -3.98%  '

# Row 11
$ws.Range("E11").Value = '  -4.97%  '

# Row 12
$ws.Range("E12").Value = '  -4.33%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '3.833.02'
$ws.Range("E13").Value = '  -5.65%  '

# Row 14
$ws.Range("E14").Value = '  -0.01%  '

# Row 15
$ws.Range("E15").Value = '  -5.61%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.267.28'
$ws.Range("E16").Value = '  -5.85%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '24.31'
$ws.Range("E17").Value = '  +0.08%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '60.400.12'
$ws.Range("E18").Value = '  -5.77%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '5.64'
$ws.Range("E19").Value = '  -0.90%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '13.26'
$ws.Range("E20").Value = '  -1.38%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '9.00'
$ws.Range("E21").Value = '  -9.83%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '352.57'
$ws.Range("E22").Value = '  -8.43%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.551'
$ws.Range("E23").Value = '  -3.46%  '

# Row 24
$ws.Range("E24").Value = '  -0.19%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '3.407.30'
$ws.Range("E25").Value = '  -5.80%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '69.38'
$ws.Range("E26").Value = '  -7.27%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.0000109'
$ws.Range("E27").Value = '  -2.15%  '

# Row 28
$ws.Range("E28").Value = '  +0.13%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.32'
$ws.Range("E29").Value = '  +3.81%  '

# Row 30
$ws.Range("E30").Value = '  -0.11%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.81'
$ws.Range("E31").Value = '  -1.38%  '

# Row 32
$ws.Range("E32").Value = '  -6.17%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("B33").Value = 'USDe'
$ws.Range("C33").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D33").Value = '1.00'
$ws.Range("E33").Value = '  +0.07%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("B34").Value = 'Kaspa'
$ws.Range("C34").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D34").Value = '0.149'
$ws.Range("E34").Value = '  -1.88%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.308.51'
$ws.Range("E35").Value = '  -5.68%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '22.68'
$ws.Range("E36").Value = '  -1.09%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.31'
$ws.Range("E37").Value = '  +1.98%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '6.80'
$ws.Range("E38").Value = '  +0.81%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.48'
$ws.Range("E39").Value = '  -0.86%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '157.85'
$ws.Range("E40").Value = '  -2.73%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0750'
$ws.Range("E41").Value = '  -3.46%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.00'
$ws.Range("E42").Value = '  +0.16%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("B43").Value = 'Filecoin'
$ws.Range("C43").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D43").Value = '4.36'
$ws.Range("E43").Value = '  +1.62%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("B44").Value = 'OKB'
$ws.Range("C44").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D44").Value = '40.89'
$ws.Range("E44").Value = '  -1.07%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.735'
$ws.Range("E45").Value = '  -7.70%  '

# Row 46
$ws.Range("E46").Value = '  +0.43%  '

# Row 47
$ws.Range("E47").Value = '  -4.30%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '22.34'
$ws.Range("E48").Value = '  -4.03%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '6.65'
$ws.Range("E49").Value = '  -0.76%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.864'
$ws.Range("E50").Value = '  -4.47%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '21.50'
$ws.Range("E51").Value = '  +5.29%  '
